# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '51.089.74'
$ws.Range("E2").Value = '  -15.88%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.226.68'
$ws.Range("E3").Value = '  -23.47%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.09%  '

# Row 5: BNB
$ws.Range("D5").Value = "'419.99"
$ws.Range("E5").Value = '  -20.53%  '

# Row 6: Solana
$ws.Range("D6").Value = "'115.25"
$ws.Range("E6").Value = '  -19.79%  '

# Row 7: USDC
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = '  -0.24%  '

# Row 8: XRP
$ws.Range("D8").Value = "'0.451"
$ws.Range("E8").Value = '  -18.76%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.215.64'
$ws.Range("E9").Value = '  -24.07%  '

# Row 10: Toncoin
$ws.Range("D10").Value = "'5.01"
$ws.Range("E10").Value = '  -16.90%  '

# Row 11: Dogecoin
$ws.Range("E11").Value = '  -21.31%  '

# Row 12: Cardano
$ws.Range("D12").Value = "'0.287"
$ws.Range("E12").Value = '  -20.42%  '

# Row 13: TRON
$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = '  -6.75%  '

# Row 14: WrappedBTC
$ws.Range("D14").Value = '51.106.99'
$ws.Range("E14").Value = '  -15.81%  '

# Row 15: Avalanche
$ws.Range("D15").Value = "'18.10"
$ws.Range("E15").Value = '  -20.74%  '

# Row 16: ShibaInu
$ws.Range("E16").Value = '  -20.83%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.228.11'
$ws.Range("E17").Value = '  -23.41%  '

# Row 18: BitcoinCash -> Polkadot
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'3.77"
$ws.Range("E18").Value = '  -25.09%  '

# Row 19: Polkadot -> BitcoinCash
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = "'285.93"
$ws.Range("E19").Value = '  -20.88%  '

# Row 20: Dai
$ws.Range("D20").Value = "'0.998"
$ws.Range("E20").Value = '  -0.17%  '

# Row 21: LEO
$ws.Range("D21").Value = "'5.58"
$ws.Range("E21").Value = '  -1.57%  '

# Row 22: Chainlink
$ws.Range("D22").Value = "'8.48"
$ws.Range("E22").Value = '  -27.66%  '

# Row 23: Uniswap
$ws.Range("D23").Value = "'5.03"
$ws.Range("E23").Value = '  -24.62%  '

# Row 24: Litecoin
$ws.Range("D24").Value = "'51.79"
$ws.Range("E24").Value = '  -20.27%  '

# Row 25: Polygon
$ws.Range("D25").Value = "'0.348"
$ws.Range("E25").Value = '  -23.59%  '

# Row 26: Kaspa
$ws.Range("D26").Value = "'0.137"
$ws.Range("E26").Value = '  -24.47%  '

# Row 27: USDe
$ws.Range("E27").Value = '  -0.38%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = "'6.58"
$ws.Range("E28").Value = '  -16.75%  '

# Row 29: PEPE
$ws.Range("D29").Value = '0.0₃0644'
$ws.Range("E29").Value = '  -24.38%  '

# Row 30: Monero
$ws.Range("D30").Value = "'139.94"
$ws.Range("E30").Value = '  -7.48%  '

# Row 31: EthereumClassic
$ws.Range("D31").Value = "'16.40"
$ws.Range("E31").Value = '  -17.39%  '

# Row 32: PancakeSwap
$ws.Range("D32").Value = "'1.26"
$ws.Range("E32").Value = '  -25.80%  '

# Row 33: Aptos
$ws.Range("D33").Value = "'4.49"
$ws.Range("E33").Value = '  -19.80%  '

# Row 34: Fetch.AI
$ws.Range("D34").Value = "'0.781"
$ws.Range("E34").Value = '  -22.45%  '

# Row 35: NEARProtocol
$ws.Range("E35").Value = '  -24.32%  '

# Row 36: FirstDigitalUSD
$ws.Range("D36").Value = "'0.996"
$ws.Range("E36").Value = '  -0.12%  '

# Row 37: ImmutableX
$ws.Range("D37").Value = "'0.951"
$ws.Range("E37").Value = '  -21.14%  '

# Row 38: WhiteBITCoin -> OKB
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = "'30.99"
$ws.Range("E38").Value = '  -18.31%  '

# Row 39: OKB -> WhiteBITCoin
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = "'10.17"
$ws.Range("E39").Value = '  -1.59%  '

# Row 40: Filecoin
$ws.Range("D40").Value = "'2.99"
$ws.Range("E40").Value = '  -19.93%  '

# Row 41: Stacks
$ws.Range("E41").Value = '  -22.52%  '

# Row 42: Hedera
$ws.Range("D42").Value = "'0.0475"
$ws.Range("E42").Value = '  -19.03%  '

# Row 43: Maker
$ws.Range("D43").Value = '1.815.82'
$ws.Range("E43").Value = '  -20.99%  '

# Row 44: Mantle
$ws.Range("D44").Value = "'0.497"
$ws.Range("E44").Value = '  -23.72%  '

# Row 45: VeChain -> Stellar
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = "'0.0796"
$ws.Range("E45").Value = '  -13.87%  '

# Row 46: Stellar -> VeChain
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = "'0.0196"
$ws.Range("E46").Value = '  -17.97%  '

# Row 47: ZEEBU -> EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'14.95"
$ws.Range("E47").Value = '  -27.24%  '

# Row 48: EnergySwap -> ZEEBU
$ws.Range("B48").Value = 'ZEEBU'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D48").Value = "'4.51"
$ws.Range("E48").Value = '  -6.87%  '

# Row 49: RenderToken
$ws.Range("D49").Value = "'3.61"
$ws.Range("E49").Value = '  -27.61%  '

# Row 50: InjectiveProtocol
$ws.Range("D50").Value = "'14.33"
$ws.Range("E50").Value = '  -22.10%  '

# Row 51: BitgetToken
$ws.Range("D51").Value = "'0.859"
$ws.Range("E51").Value = '  -15.45%  '
